$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 532, pushing existing rows 532:600 down by one.
$ws.Rows.Item(532).Insert()

# Populate the new row 532 with the new weekly record.
$ws.Cells.Item(532, 1).Value = 6
$ws.Cells.Item(532, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(532, 3).Value = "Metropolitana"
$ws.Cells.Item(532, 4).Value = 45127
$ws.Cells.Item(532, 5).Value = 13
$ws.Cells.Item(532, 6).Value = 100112032
$ws.Cells.Item(532, 7).Value = "Zapallo italiano"
$ws.Cells.Item(532, 8).Value = "Sin especificar"
$ws.Cells.Item(532, 9).Value = "Primera"
$ws.Cells.Item(532, 10).Value = 270
$ws.Cells.Item(532, 11).Value = 13000
$ws.Cells.Item(532, 12).Value = 14000
$ws.Cells.Item(532, 13).Value = 13556
$ws.Cells.Item(532, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(532, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(532, 16).Value = 271
$ws.Cells.Item(532, 17).Value = 50
$ws.Cells.Item(532, 18).Value = "Hortaliza"
